# NRC Reactor Status Download and Read CSV
# Inserts two new date columns (G:H) pushing the previous G/H data to I/J,
# populates the new column with the latest NRC capacity-factor date,
# refreshes a couple of capacity figures, adds a total-capacity formula,
# a footnote about Vogtle 3 starting commercial operation, and a small
# scratch area (rows 17-19) with source link + capacity-list math.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at G:H - this shifts the existing G,H columns
# (and all their formulas/references) over to I,J automatically.
$ws.Columns("G:H").Insert()

# New G column holds the newest NRC report date; match the date format
# used by the neighboring (now shifted) date columns.
$ws.Range("G1").Value = 45214
$ws.Range("G1").NumberFormat = $ws.Range("I1").NumberFormat

# Updated Farley 2 capacity figure.
$ws.Range("C3").Value = 896

# Footnote next to Vogtle 3's row.
$ws.Range("L8").Value = "Began commercial operation 8/31"

# New total nameplate capacity formula.
$ws.Range("C11").Formula = "=SUM(C2:C8)"

# New scratch area: source link + a quick capacity-list calculation.
$ws.Range("C17").Value = "https://www.eia.gov/nuclear/reactors/reactorcapacity.php"
$ws.Range("A17").Value = "Capacity List"
$ws.Range("O17").Formula = "=AVERAGE(3.22, 3.75)"
$ws.Range("Q17").Formula = "=AVERAGE(3.22, 3.86)"
$ws.Range("O18").Formula = "=O17-2.35"
$ws.Range("O19").Formula = "=O17-2.45"

# Restore the active selection to match the saved view.
$ws.Range("F5").Select() | Out-Null
